$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010-18")
$ws.Activate()

# --- New row 24: a fresh "Baseline 2010-18 C445" data row copied/adjusted from row 14 ---
$ws.Range("A24").Value = "CW3M"
$ws.Range("B24").Value = "Baseline 2010-18 C445"
$ws.Range("C24").Value = "2010-18"

$ws.Range("D24").Value = 538.12494244444451
$ws.Range("E24").Value = 2094.2995878888887
$ws.Range("F24").Value = 5.8220211111111109
$ws.Range("G24").Value = 332.2750817777777
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 8.1971097777777775
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 596.67409588888893
$ws.Range("L24").Value = 91.777595333333338
$ws.Range("M24").Value = 1765.9676106666668
$ws.Range("N24").Value = 524.34534722222224
$ws.Range("O24").Value = 15018.345052222223
$ws.Range("P24").Value = 2216.8192002222222
$ws.Range("Q24").Value = 0.045906222222222182
$ws.Range("R24").Value = -0.0000089999999999999901

$ws.Range("D24:N24").NumberFormat = "0.00"
$ws.Range("O24").NumberFormat = "0"
$ws.Range("O24").Interior.Color = 65535
$ws.Range("P24").NumberFormat = "0"
$ws.Range("Q24").NumberFormat = "0.00"
$ws.Range("R24").NumberFormat = "0.000000"

# --- Selection / scroll state left after the edit: whole row 24 selected, view scrolled to top ---
$ws.Range("A24:XFD24").Select()
